$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Assurance MM6" milestone label (and its derived cells) is being
# renamed to "Assurance MM1" throughout column A.
for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value2
    if ($v -ne $null -and $v -like "*Assurance MM6*") {
        $cell.Value2 = $v -replace "Assurance MM6", "Assurance MM1"
    }
}

# Move the sheet's active-cell selection from B40 to A9.
$ws.Range("A9").Select()
